$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column C flips from "No" (red fill, style carried from C43/etc.)
# to "Yes" (green fill, same style as the other completed rows, e.g. C6).
# Use copy/paste-special of formats so the existing "Yes" style (fill/font)
# is reused instead of minting a brand new style entry.
$yesRows        = @(43, 46, 55, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76)
# Of those, these additionally gain a completion date in column D (matching
# the date style already used elsewhere, e.g. D56).
$dateRows       = @(55, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76)
$completionDate = 43635

foreach ($r in $yesRows) {
    $ws.Range("C6").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("C$r").Value = "Yes"
}

foreach ($r in $dateRows) {
    $ws.Range("D56").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("D$r").Value = $completionDate
}

# Move the active selection from B58 to C46 (test 2 finished, test 3 starting).
$ws.Range("C46").Select() | Out-Null
